# Generate Report for Archive
#
# 1. The localization status moves from "Ready for handoff" to
#    "In Translation" everywhere that status string appears (Overview
#    sheet's zh-cn/de-de columns, plus the Status column on each
#    per-locale sheet).
# 2. The (now narrower) Status-ish columns are resized down to match.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

$newStatus = "In Translation"

# --- Overview sheet: columns E (zh-cn) and F (de-de), rows 2-3 ---
$ws1.Range("E2").Value = $newStatus
$ws1.Range("F2").Value = $newStatus
$ws1.Range("E3").Value = $newStatus
$ws1.Range("F3").Value = $newStatus

# --- zh-cn / de-de sheets: Status column C, rows 2-3 ---
$ws2.Range("C2").Value = $newStatus
$ws2.Range("C3").Value = $newStatus
$ws3.Range("C2").Value = $newStatus
$ws3.Range("C3").Value = $newStatus

# --- Shrink the columns that held the longer status text ---
$ws1.Columns.Item(5).ColumnWidth = 12.5   # E (zh-cn)
$ws1.Columns.Item(6).ColumnWidth = 12.5   # F (de-de)
$ws2.Columns.Item(3).ColumnWidth = 12.5   # C (Status)
$ws3.Columns.Item(3).ColumnWidth = 12.5   # C (Status)
